$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks numeric but must stay stored as text
# (matches the original inlineStr formatting, e.g. "7.50", "1.60", "0.989").
# Force a text number-format first so Excel does not coerce them to numbers
# (which would also silently drop significant trailing zeros).
$riskyRefs = @('D4','D5','D6','D7','D8','D10','D14','D15','D17','D18','D19','D21','D22','D23','D24','D25','D26','D27','D29','D31','D32','D34','D35','D36','D39','D40','D41','D42','D43','D44','D45','D46','D48','D49','D51')
foreach ($r in $riskyRefs) { $ws.Range($r).NumberFormat = "@" }

$ws.Range('D2').Value = '27.601.98'
$ws.Range('E2').Value = '  +2.10%  '
$ws.Range('D3').Value = '1.568.65'
$ws.Range('E3').Value = '  +0.46%  '
$ws.Range('D4').Value = '0.989'
$ws.Range('E4').Value = '  -1.54%  '
$ws.Range('D5').Value = '210.94'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('D6').Value = '0.493'
$ws.Range('E6').Value = '  +0.55%  '
$ws.Range('D7').Value = '0.988'
$ws.Range('E7').Value = '  -1.73%  '
$ws.Range('D8').Value = '23.19'
$ws.Range('E8').Value = '  +5.74%  '
$ws.Range('E9').Value = '  +1.40%  '
$ws.Range('D10').Value = '0.0597'
$ws.Range('E10').Value = '  +0.37%  '
$ws.Range('E11').Value = '  +1.66%  '
$ws.Range('D12').Value = '1.796.32'
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').Value = '1.612.76'
$ws.Range('E13').Value = '  +3.30%  '
$ws.Range('D14').Value = '3.76'
$ws.Range('E14').Value = '  -0.14%  '
$ws.Range('D15').Value = '0.521'
$ws.Range('E15').Value = '  +0.67%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '27.553.43'
$ws.Range('E16').Value = '  +1.94%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').Value = '63.07'
$ws.Range('E17').Value = '  +1.86%  '
$ws.Range('D18').Value = '229.34'
$ws.Range('E18').Value = '  +6.59%  '
$ws.Range('D19').Value = '7.50'
$ws.Range('E19').Value = '  +1.46%  '
$ws.Range('D20').Value = '0.0₃0704'
$ws.Range('E20').Value = '  +0.24%  '
$ws.Range('D21').Value = '0.989'
$ws.Range('E21').Value = '  -1.63%  '
$ws.Range('D22').Value = '4.11'
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('D23').Value = '9.42'
$ws.Range('E23').Value = '  +2.96%  '
$ws.Range('D24').Value = '1.95'
$ws.Range('E24').Value = '  +0.93%  '
$ws.Range('D25').Value = '149.69'
$ws.Range('E25').Value = '  -2.51%  '
$ws.Range('D26').Value = '15.29'
$ws.Range('E26').Value = '  +1.53%  '
$ws.Range('D27').Value = '6.59'
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('E28').Value = '  +1.68%  '
$ws.Range('D29').Value = '0.990'
$ws.Range('E29').Value = '  -1.38%  '
$ws.Range('E30').Value = '  +0.53%  '
$ws.Range('D31').Value = '0.0473'
$ws.Range('E31').Value = '  +0.35%  '
$ws.Range('D32').Value = '3.24'
$ws.Range('E32').Value = '  +0.56%  '
$ws.Range('D33').Value = '1.456.37'
$ws.Range('E33').Value = '  +1.92%  '
$ws.Range('D34').Value = '3.12'
$ws.Range('E34').Value = '  -1.47%  '
$ws.Range('D35').Value = '1.60'
$ws.Range('E35').Value = '  -0.53%  '
$ws.Range('D36').Value = '1.05'
$ws.Range('E36').Value = '  -5.97%  '
$ws.Range('E37').Value = '  -0.52%  '
$ws.Range('E38').Value = '  +0.80%  '
$ws.Range('D39').Value = '0.542'
$ws.Range('E39').Value = '  +2.12%  '
$ws.Range('D40').Value = '0.812'
$ws.Range('E40').Value = '  +0.71%  '
$ws.Range('D41').Value = '2.39'
$ws.Range('E41').Value = '  -2.09%  '
$ws.Range('D42').Value = '5.66'
$ws.Range('E42').Value = '  -2.92%  '
$ws.Range('D43').Value = '0.989'
$ws.Range('D44').Value = '1.85'
$ws.Range('E44').Value = '  +6.50%  '
$ws.Range('D45').Value = '0.972'
$ws.Range('E45').Value = '  -2.60%  '
$ws.Range('D46').Value = '64.03'
$ws.Range('E46').Value = '  -0.36%  '
$ws.Range('D47').Value = '1.706.30'
$ws.Range('E47').Value = '  +0.18%  '
$ws.Range('D48').Value = '86.84'
$ws.Range('E48').Value = '  +1.54%  '
$ws.Range('D49').Value = '0.0523'
$ws.Range('E49').Value = '  +1.29%  '
$ws.Range('D50').Value = '0.0₇0977'
$ws.Range('E50').Value = '  -3.81%  '
$ws.Range('D51').Value = '39.70'
$ws.Range('E51').Value = '  +17.35%  '
